$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": add a handoff-report row for the new file that just
# finished handback (da7bd816-...md), mirroring the existing rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A7").Value = "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
$wsOverview.Range("B7").Value = "e2e\da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-07 08:34:01"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a8b1191b915c6ca65b09d53d5270eb4c86bf375/e2e/da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md",
    "",
    "",
    "e2e\da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": add the matching handoff row to the zh-cn status table.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A7").Value = "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.7a000bc0a9c0a53acc837595595b64f67a126d72.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-07 08:33:50"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a8b1191b915c6ca65b09d53d5270eb4c86bf375/e2e/da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md",
    "",
    "",
    "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": add the matching handoff row to the de-de status table.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A7").Value = "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.7a000bc0a9c0a53acc837595595b64f67a126d72.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-07 08:34:01"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A7"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a8b1191b915c6ca65b09d53d5270eb4c86bf375/e2e/da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md",
    "",
    "",
    "da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
) | Out-Null

Write-Host "Handoff report row added for da7bd816-6caf-4b1f-b1dd-aaabfa97f568.md"
